$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Append the new log row (row 7)
$ws.Range("A7").Value = "Kun jij dit even regelen?"
$ws.Range("B7").Value = "mailmind.test@zohomail.eu"
$ws.Range("C7").Value = "Testmail #1: Kun jij dit even regelen?"
$ws.Range("D7").Value = "Overig"
$ws.Range("E7").Value = "Geachte klant,`nBedankt voor uw e-mail. Om u beter van dienst te kunnen zijn, zou u wat meer informatie kunnen verstrekken over wat u precies wilt regelen?`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$ws.Range("F7").Value = "2025-07-27 17:54:51"
$ws.Range("G7").Value = "Ja"
$ws.Range("H7").Value = "Nee"
$ws.Range("I7").Value = "Ja"
$ws.Range("J7").Value = "Ja"

# Keep the new row's height consistent with the rest of the sheet (writing the
# multi-line E7 value would otherwise trigger an automatic row resize)
$ws.Rows.Item(7).RowHeight = $ws.Rows.Item(6).RowHeight

# Extend the conditional formatting ranges to cover the new row
foreach ($col in @("D", "G", "H", "I", "J")) {
    $fcs = $ws.Range("$col`2:$col`6").FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($ws.Range("$col`2:$col`7"))
    }
}

# Update the Dashboard aggregate count
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 6
